$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.782.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.12%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.074.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.46%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.73'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.78%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.540'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.075.20'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.37%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.92'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.11%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.21%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.121'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.581.57'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.20'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.663.35'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.073.79'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '484.00'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.61%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.81%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.51%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.46'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.05'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.93'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.65'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +7.42%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.14%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.52%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.32'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.21%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0823'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.23%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.84%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.25'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.80%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '50.77'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("B41").Value = 'Cosmos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.22'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '440.56'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.92%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.114'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +4.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0365'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.43%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.844.07'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.10'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.84'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.27%  '

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.64%  '
